$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''64.225.95'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -3.18%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''3.141.61'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -2.31%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.02%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''607.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -0.26%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''147.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -5.42%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = '''  +0.16%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''3.140.27'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  -2.38%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.529'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -3.81%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = '''  -5.14%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''5.58'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -2.21%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''0.475'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -5.33%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''0.0000259'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Value = '''36.64'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -4.72%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''3.655.85'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -2.31%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''64.231.54'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -3.29%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = '''  +0.01%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''3.147.04'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -2.13%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''6.97'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  -4.23%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''480.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -5.22%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''14.60'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -4.49%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''0.709'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -3.03%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = '''  -3.83%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''13.77'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -5.59%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''83.84'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -1.50%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  +0.09%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''2.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -2.33%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = '''  -5.42%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = '''  -4.94%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''0.124'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  -9.23%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''6.87'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -1.03%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = '''Stacks'
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = '''https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = '''2.73'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -5.55%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value = '''FirstDigitalUSD'
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = '''https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = '''0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -0.26%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''26.70'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -5.54%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '''  -5.34%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''6.09'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -5.44%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''54.46'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -1.69%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''3.15'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  +3.43%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''0.0₃0737'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  -4.93%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''451.98'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -10.06%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''0.0399'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -4.78%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = '''  -5.87%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = '''  -3.31%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''2.874.20'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -1.71%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = '''  -8.44%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''2.27'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -6.65%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = '''  -5.29%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = '''  -0.03%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = '''  -3.51%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''0.115'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -2.40%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''120.25'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -1.36%  '
$ws.Range('E51').Style = 'Normal'
